# Scheduled market-data refresh: re-pull current Market Board prices
# (currentAveragePrice / NQ / HQ) and recompute Leve crafting profit
# columns (LevePriceNQ/HQ, LeveProfitNQ/HQ) for every job sheet.
$wb = $excel.ActiveWorkbook

# --- ALC sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 8999.799999999999
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10349

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 25813.857
$ws.Range("J43").Value = 14974.5
$ws.Range("L43").Value = 14974.5
$ws.Range("N43").Value = -15112.5

# Row 46: Always Have an Exit Plan
$ws.Range("H46").Value = 8873.333000000001
$ws.Range("J46").Value = 8873.333000000001
$ws.Range("L46").Value = 26619.999
$ws.Range("N46").Value = -26857.999

# Row 60: Make Up Your Mind or Else
$ws.Range("H60").Value = 8873.333000000001
$ws.Range("J60").Value = 8873.333000000001
$ws.Range("L60").Value = 26619.999
$ws.Range("N60").Value = -27587.999

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 4056.5557
$ws.Range("I62").Value = 4081.08
$ws.Range("K62").Value = 4081.08
$ws.Range("M62").Value = -3457.08

# Row 64: Forged from the Void
$ws.Range("H64").Value = 5768.522
$ws.Range("I64").Value = 3485.5334
$ws.Range("J64").Value = 10049.125
$ws.Range("K64").Value = 3485.5334
$ws.Range("L64").Value = 10049.125
$ws.Range("M64").Value = -3237.5334
$ws.Range("N64").Value = -10545.125

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 4056.5557
$ws.Range("I65").Value = 4081.08
$ws.Range("K65").Value = 20405.4
$ws.Range("M65").Value = -17285.4

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 5768.522
$ws.Range("I67").Value = 3485.5334
$ws.Range("J67").Value = 10049.125
$ws.Range("K67").Value = 3485.5334
$ws.Range("L67").Value = 10049.125
$ws.Range("M67").Value = -2627.5334
$ws.Range("N67").Value = -11765.125

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 4073.3635
$ws.Range("J86").Value = 4363.4
$ws.Range("L86").Value = 4363.4
$ws.Range("N86").Value = -6609.4

# Row 87: There Was a Late Fee
$ws.Range("H87").Value = 43936.582
$ws.Range("J87").Value = 43936.582
$ws.Range("L87").Value = 43936.582
$ws.Range("N87").Value = -46432.582

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 4073.3635
$ws.Range("J89").Value = 4363.4
$ws.Range("L89").Value = 21817
$ws.Range("N89").Value = -33049

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("H90").Value = 43936.582
$ws.Range("J90").Value = 43936.582
$ws.Range("L90").Value = 131809.746
$ws.Range("N90").Value = -144289.746


# --- ARM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 5987.338
$ws.Range("I32").Value = 3572.7942
$ws.Range("K32").Value = 3572.7942
$ws.Range("M32").Value = -3285.7942

# Row 43: They've Got Legs
$ws.Range("H43").Value = 45000
$ws.Range("I43").Value = 45000
$ws.Range("K43").Value = 45000
$ws.Range("M43").Value = -44687

# Row 69: The Cut Alembical Cord
$ws.Range("H69").Value = 496845
$ws.Range("J69").Value = 496845
$ws.Range("L69").Value = 496845
$ws.Range("N69").Value = -498343

# Row 72: Sheer Distill Power (L)
$ws.Range("H72").Value = 496845
$ws.Range("J72").Value = 496845
$ws.Range("L72").Value = 1490535
$ws.Range("N72").Value = -1498023

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1801.7142
$ws.Range("I122").Value = 1490
$ws.Range("J122").Value = 2362.8
$ws.Range("K122").Value = 4470
$ws.Range("L122").Value = 7088.400000000001
$ws.Range("M122").Value = -2020
$ws.Range("N122").Value = -11988.4

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2813.1892
$ws.Range("I132").Value = 2640.9688
$ws.Range("K132").Value = 7922.9064
$ws.Range("M132").Value = -5392.9064


# --- BSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 27: Claw Daddy
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2477.8
$ws.Range("I86").Value = 2863
$ws.Range("J86").Value = 1900
$ws.Range("K86").Value = 2863
$ws.Range("L86").Value = 1900
$ws.Range("M86").Value = -1740
$ws.Range("N86").Value = -4146

# Row 88: Swords for Plowshares
$ws.Range("H88").Value = 129620
$ws.Range("J88").Value = 129620
$ws.Range("L88").Value = 129620
$ws.Range("N88").Value = -130432

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2477.8
$ws.Range("I89").Value = 2863
$ws.Range("J89").Value = 1900
$ws.Range("K89").Value = 14315
$ws.Range("L89").Value = 9500
$ws.Range("M89").Value = -8699
$ws.Range("N89").Value = -20732

# Row 91: Negative, They Are Meat Popsicles (L)
$ws.Range("H91").Value = 129620
$ws.Range("J91").Value = 129620
$ws.Range("L91").Value = 129620
$ws.Range("N91").Value = -132428

# Row 94: High Steal
$ws.Range("H94").Value = 2045.421
$ws.Range("I94").Value = 1246.8572
$ws.Range("K94").Value = 1246.8572
$ws.Range("M94").Value = -795.8571999999999

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 14002.75
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 14002.75
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 14002.75
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -16998.75


# --- CRP sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 394.4
$ws.Range("I22").Value = 368.125
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 368.125
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = -18.125
$ws.Range("N22").Value = -1199.5

# Row 31: Wall Not Found
$ws.Range("H31").Value = 39182.406
$ws.Range("I31").Value = 43609.375
$ws.Range("K31").Value = 43609.375
$ws.Range("M31").Value = -43314.375

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 39182.406
$ws.Range("I34").Value = 43609.375
$ws.Range("K34").Value = 43609.375
$ws.Range("M34").Value = -43407.375

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2846.2632
$ws.Range("I58").Value = 2767.6785
$ws.Range("J58").Value = 3066.3
$ws.Range("K58").Value = 2767.6785
$ws.Range("L58").Value = 3066.3
$ws.Range("M58").Value = -2564.6785
$ws.Range("N58").Value = -3472.3

# Row 86: Birch, Please
$ws.Range("H86").Value = 8416
$ws.Range("I86").Value = 9669
$ws.Range("J86").Value = 7998.3335
$ws.Range("K86").Value = 9669
$ws.Range("L86").Value = 7998.3335
$ws.Range("M86").Value = -8546
$ws.Range("N86").Value = -10244.3335

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 8416
$ws.Range("I89").Value = 9669
$ws.Range("J89").Value = 7998.3335
$ws.Range("K89").Value = 48345
$ws.Range("L89").Value = 39991.6675
$ws.Range("M89").Value = -42729
$ws.Range("N89").Value = -51223.6675

# Row 107: Built to Last
$ws.Range("H107").Value = 1012.05
$ws.Range("I107").Value = 679.5
$ws.Range("J107").Value = 1154.5714
$ws.Range("K107").Value = 679.5
$ws.Range("L107").Value = 1154.5714
$ws.Range("M107").Value = 1240.5
$ws.Range("N107").Value = -4994.5714

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3079.4243
$ws.Range("I132").Value = 2887.92
$ws.Range("J132").Value = 3677.875
$ws.Range("K132").Value = 8663.76
$ws.Range("L132").Value = 11033.625
$ws.Range("M132").Value = -6133.76
$ws.Range("N132").Value = -16093.625

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 15613.706
$ws.Range("I134").Value = 5312.2334
$ws.Range("K134").Value = 15936.7002
$ws.Range("M134").Value = -13401.7002

# Row 136: Turali Quality
$ws.Range("H136").Value = 2846.2632
$ws.Range("I136").Value = 2767.6785
$ws.Range("J136").Value = 3066.3
$ws.Range("K136").Value = 8303.0355
$ws.Range("L136").Value = 9198.900000000001
$ws.Range("M136").Value = -5753.0355
$ws.Range("N136").Value = -14298.9

# Row 137: Lament of the Lazylump
$ws.Range("H137").Value = 87473.336
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 87473.336
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 87473.336
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -97673.336


# --- CUL sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 8: Whip It
$ws.Range("H8").Value = 4899
$ws.Range("I8").Value = 4899
$ws.Range("K8").Value = 14697
$ws.Range("M8").Value = -14558

# Row 132: More Mezcal
$ws.Range("H132").Value = 1297.3636
$ws.Range("I132").Value = 1318.0625
$ws.Range("J132").Value = 1242.1666
$ws.Range("K132").Value = 11862.5625
$ws.Range("L132").Value = 11179.4994
$ws.Range("M132").Value = -9332.5625
$ws.Range("N132").Value = -16239.4994


# --- GSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2779.1052
$ws.Range("I80").Value = 2708.077
$ws.Range("K80").Value = 2708.077
$ws.Range("M80").Value = -1710.077

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2779.1052
$ws.Range("I83").Value = 2708.077
$ws.Range("K83").Value = 13540.385
$ws.Range("M83").Value = -8548.385000000002

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2060.926
$ws.Range("I113").Value = 2275.238
$ws.Range("J113").Value = 1310.8334
$ws.Range("K113").Value = 2275.238
$ws.Range("L113").Value = 1310.8334
$ws.Range("M113").Value = -105.2379999999998
$ws.Range("N113").Value = -5650.8334

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2647.3696
$ws.Range("I122").Value = 2245.7693
$ws.Range("J122").Value = 4884.857
$ws.Range("K122").Value = 6737.3079
$ws.Range("L122").Value = 14654.571
$ws.Range("M122").Value = -4287.3079
$ws.Range("N122").Value = -19554.571

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2798.5625
$ws.Range("I132").Value = 2354.262
$ws.Range("K132").Value = 7062.786
$ws.Range("M132").Value = -4532.786


# --- LTW sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 16: Saddle Sore
$ws.Range("H16").Value = 2389.48
$ws.Range("I16").Value = 2467.1738
$ws.Range("K16").Value = 2467.1738
$ws.Range("M16").Value = -2297.1738

# Row 40: Best Served Toad
$ws.Range("H40").Value = 3142.7334
$ws.Range("I40").Value = 2387.32
$ws.Range("J40").Value = 6919.8
$ws.Range("K40").Value = 2387.32
$ws.Range("L40").Value = 6919.8
$ws.Range("M40").Value = -2251.32
$ws.Range("N40").Value = -7191.8

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2316.6667
$ws.Range("I46").Value = 1430
$ws.Range("J46").Value = 3425
$ws.Range("K46").Value = 1430
$ws.Range("L46").Value = 3425
$ws.Range("M46").Value = -1242
$ws.Range("N46").Value = -3801

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 20124.5
$ws.Range("I100").Value = 10166.333
$ws.Range("K100").Value = 10166.333
$ws.Range("M100").Value = -9625.333000000001


# --- WVR sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1992.4348
$ws.Range("I126").Value = 1851.4375
$ws.Range("K126").Value = 5554.3125
$ws.Range("M126").Value = -3084.3125

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1723.7059
$ws.Range("I132").Value = 1634.0646
$ws.Range("K132").Value = 4902.1938
$ws.Range("M132").Value = -2372.1938


Write-Output "Applied all cell updates"